# Added a missing interval in schedule template test file
# The schedule table is missing the "16:00 - 17:00" slot between
# "15:00 - 16:00" (row 9) and "17:00 - 18:00" (row 10, before the edit).
# Insert a new row at position 10, push the remaining rows down, and fill
# in the missing interval.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 10 ("17:00 - 18:00") and row 11 ("18:00 - 19:00") down by one
# row to make room for the new interval.
$ws.Rows.Item(10).Insert()

$ws.Range("A10").Value = "16:00 - 17:00"
$ws.Range("A11").Value = "17:00 - 18:00"
$ws.Range("A12").Value = "18:00 - 19:00"

# Match the saved selection/active cell from the authored workbook.
$ws.Range("A10").Select()

# Restore default print setup (A4 portrait) as captured by the edit.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
